$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.852.37"
$ws.Range("D3").Value = "1.735.89"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2760"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "1.737.89"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07079"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6417"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "25.839.28"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006659"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "1.959.63"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.243"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.801"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.513"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.786"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08322"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.716"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.530"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04491"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9730"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6195"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.675"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01574"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.911"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3854"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.030"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05321"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.209"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.646"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
